$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-27 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-28 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("895÷2=447, 1", $true, $false, $false, $false, $false, $true, 1, $false, "231÷6=38, 3", 2) | Out-Null
$d.Content.Find.Execute("214÷4=53, 2", $true, $false, $false, $false, $false, $true, 1, $false, "331÷3=110, 1", 2) | Out-Null
$d.Content.Find.Execute("158÷4=39, 2", $true, $false, $false, $false, $false, $true, 1, $false, "872÷3=290, 2", 2) | Out-Null
$d.Content.Find.Execute("429÷2=214, 1", $true, $false, $false, $false, $false, $true, 1, $false, "690÷9=76, 6", 2) | Out-Null
$d.Content.Find.Execute("736÷6=122, 4", $true, $false, $false, $false, $false, $true, 1, $false, "791÷2=395, 1", 2) | Out-Null
$d.Content.Find.Execute("215÷2=107, 1", $true, $false, $false, $false, $false, $true, 1, $false, "934÷2=467, 0", 2) | Out-Null
$d.Content.Find.Execute("411÷9=45, 6", $true, $false, $false, $false, $false, $true, 1, $false, "718÷4=179, 2", 2) | Out-Null
$d.Content.Find.Execute("834÷8=104, 2", $true, $false, $false, $false, $false, $true, 1, $false, "817÷7=116, 5", 2) | Out-Null
$d.Content.Find.Execute("999÷6=166, 3", $true, $false, $false, $false, $false, $true, 1, $false, "334÷6=55, 4", 2) | Out-Null
$d.Content.Find.Execute("515÷3=171, 2", $true, $false, $false, $false, $false, $true, 1, $false, "732÷6=122, 0", 2) | Out-Null
$d.Content.Find.Execute("269÷4=67, 1", $true, $false, $false, $false, $false, $true, 1, $false, "475÷2=237, 1", 2) | Out-Null
$d.Content.Find.Execute("564÷8=70, 4", $true, $false, $false, $false, $false, $true, 1, $false, "458÷4=114, 2", 2) | Out-Null
$d.Content.Find.Execute("103÷6=17, 1", $true, $false, $false, $false, $false, $true, 1, $false, "433÷3=144, 1", 2) | Out-Null
$d.Content.Find.Execute("682÷8=85, 2", $true, $false, $false, $false, $false, $true, 1, $false, "859÷4=214, 3", 2) | Out-Null
$d.Content.Find.Execute("149÷8=18, 5", $true, $false, $false, $false, $false, $true, 1, $false, "370÷9=41, 1", 2) | Out-Null
$d.Content.Find.Execute("879÷6=146, 3", $true, $false, $false, $false, $false, $true, 1, $false, "683÷6=113, 5", 2) | Out-Null
$d.Content.Find.Execute("253÷8=31, 5", $true, $false, $false, $false, $false, $true, 1, $false, "801÷5=160, 1", 2) | Out-Null
$d.Content.Find.Execute("106÷8=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "439÷5=87, 4", 2) | Out-Null
$d.Content.Find.Execute("511÷9=56, 7", $true, $false, $false, $false, $false, $true, 1, $false, "941÷2=470, 1", 2) | Out-Null
$d.Content.Find.Execute("554÷5=110, 4", $true, $false, $false, $false, $false, $true, 1, $false, "742÷4=185, 2", 2) | Out-Null
$d.Content.Find.Execute("583÷9=64, 7", $true, $false, $false, $false, $false, $true, 1, $false, "480÷5=96, 0", 2) | Out-Null
$d.Content.Find.Execute("614÷2=307, 0", $true, $false, $false, $false, $false, $true, 1, $false, "145÷5=29, 0", 2) | Out-Null
$d.Content.Find.Execute("606÷8=75, 6", $true, $false, $false, $false, $false, $true, 1, $false, "883÷9=98, 1", 2) | Out-Null
$d.Content.Find.Execute("170÷9=18, 8", $true, $false, $false, $false, $false, $true, 1, $false, "577÷8=72, 1", 2) | Out-Null
$d.Content.Find.Execute("919÷5=183, 4", $true, $false, $false, $false, $false, $true, 1, $false, "425÷4=106, 1", 2) | Out-Null
